$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value2 = 30000
$ws.Range("J21").Value2 = 30000
$ws.Range("L21").Value2 = 30000
$ws.Range("N21").Value2 = -30936
$ws.Range("H23").Value2 = 30000
$ws.Range("J23").Value2 = 30000
$ws.Range("L23").Value2 = 30000
$ws.Range("N23").Value2 = -30468
$ws.Range("H40").Value2 = 15008.883
$ws.Range("I40").Value2 = 4525.1665
$ws.Range("K40").Value2 = 4525.1665
$ws.Range("M40").Value2 = -4350.1665
$ws.Range("H64").Value2 = 8197.267
$ws.Range("I64").Value2 = 7198.4
$ws.Range("K64").Value2 = 7198.4
$ws.Range("M64").Value2 = -6950.4
$ws.Range("H67").Value2 = 8197.267
$ws.Range("I67").Value2 = 7198.4
$ws.Range("K67").Value2 = 7198.4
$ws.Range("M67").Value2 = -6340.4
$ws.Range("H70").Value2 = 115547.664
$ws.Range("J70").Value2 = 171239.67
$ws.Range("L70").Value2 = 513719.01
$ws.Range("N70").Value2 = -514259.01
$ws.Range("H73").Value2 = 115547.664
$ws.Range("J73").Value2 = 171239.67
$ws.Range("L73").Value2 = 513719.01
$ws.Range("N73").Value2 = -515591.01
$ws.Range("H74").Value2 = 3922.25
$ws.Range("I74").Value2 = 3922.25
$ws.Range("K74").Value2 = 3922.25
$ws.Range("M74").Value2 = -2986.25
$ws.Range("H77").Value2 = 3922.25
$ws.Range("I77").Value2 = 3922.25
$ws.Range("K77").Value2 = 19611.25
$ws.Range("M77").Value2 = -14931.25
$ws.Range("H116").Value2 = 5100.077
$ws.Range("I116").Value2 = 5210.5557
$ws.Range("J116").Value2 = 4851.5
$ws.Range("K116").Value2 = 5210.5557
$ws.Range("L116").Value2 = 4851.5
$ws.Range("M116").Value2 = -1768.5557
$ws.Range("N116").Value2 = -11735.5
$ws.Range("H137").Value2 = 1810.52
$ws.Range("I137").Value2 = 1368.5883
$ws.Range("K137").Value2 = 4105.7649
$ws.Range("M137").Value2 = -1555.7649

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2546.12
$ws.Range("I2").Value2 = 2024.9375
$ws.Range("J2").Value2 = 3472.6667
$ws.Range("K2").Value2 = 2024.9375
$ws.Range("L2").Value2 = 3472.6667
$ws.Range("M2").Value2 = -1911.9375
$ws.Range("N2").Value2 = -3698.6667
$ws.Range("H32").Value2 = 7183.7817
$ws.Range("I32").Value2 = 741.25
$ws.Range("K32").Value2 = 741.25
$ws.Range("M32").Value2 = -454.25
$ws.Range("H45").Value2 = 14998
$ws.Range("I45").Value2 = 22446.4
$ws.Range("K45").Value2 = 22446.4
$ws.Range("M45").Value2 = -22069.4
$ws.Range("H74").Value2 = 4486.778
$ws.Range("I74").Value2 = 3437.6
$ws.Range("J74").Value2 = 5798.25
$ws.Range("K74").Value2 = 3437.6
$ws.Range("L74").Value2 = 5798.25
$ws.Range("M74").Value2 = -2563.6
$ws.Range("N74").Value2 = -7546.25
$ws.Range("H77").Value2 = 4486.778
$ws.Range("I77").Value2 = 3437.6
$ws.Range("J77").Value2 = 5798.25
$ws.Range("K77").Value2 = 17188
$ws.Range("L77").Value2 = 28991.25
$ws.Range("M77").Value2 = -12820
$ws.Range("N77").Value2 = -37727.25
$ws.Range("H114").Value2 = 31132.666
$ws.Range("J114").Value2 = 31132.666
$ws.Range("L114").Value2 = 31132.666
$ws.Range("N114").Value2 = -39810.666
$ws.Range("H116").Value2 = 2546.12
$ws.Range("I116").Value2 = 2024.9375
$ws.Range("J116").Value2 = 3472.6667
$ws.Range("K116").Value2 = 2024.9375
$ws.Range("L116").Value2 = 3472.6667
$ws.Range("M116").Value2 = 269.0625
$ws.Range("N116").Value2 = -8060.6667
$ws.Range("H132").Value2 = 6882.875
$ws.Range("I132").Value2 = 6827.2856
$ws.Range("K132").Value2 = 20481.8568
$ws.Range("M132").Value2 = -17951.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2546.12
$ws.Range("I3").Value2 = 2024.9375
$ws.Range("J3").Value2 = 3472.6667
$ws.Range("K3").Value2 = 2024.9375
$ws.Range("L3").Value2 = 3472.6667
$ws.Range("M3").Value2 = -1910.9375
$ws.Range("N3").Value2 = -3700.6667
$ws.Range("H86").Value2 = 5373
$ws.Range("I86").Value2 = 4533
$ws.Range("K86").Value2 = 4533
$ws.Range("M86").Value2 = -3410
$ws.Range("H89").Value2 = 5373
$ws.Range("I89").Value2 = 4533
$ws.Range("K89").Value2 = 22665
$ws.Range("M89").Value2 = -17049
$ws.Range("H107").Value2 = 3290
$ws.Range("I107").Value2 = 3290
$ws.Range("J107").Value2 = 0
$ws.Range("K107").Value2 = 3290
$ws.Range("L107").Value2 = 0
$ws.Range("M107").Value2 = -1370
$ws.Range("N107").ClearContents()
$ws.Range("H138").Value2 = 94999
$ws.Range("J138").Value2 = 94999
$ws.Range("L138").Value2 = 94999
$ws.Range("N138").Value2 = -105279

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value2 = 45174.668
$ws.Range("I43").Value2 = 0
$ws.Range("J43").Value2 = 45174.668
$ws.Range("K43").Value2 = 0
$ws.Range("L43").Value2 = 45174.668
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value2 = -45542.668
$ws.Range("H101").Value2 = 45174.668
$ws.Range("I101").Value2 = 0
$ws.Range("J101").Value2 = 45174.668
$ws.Range("K101").Value2 = 0
$ws.Range("L101").Value2 = 45174.668
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value2 = -51664.668
$ws.Range("H105").Value2 = 1193.3636
$ws.Range("I105").Value2 = 1241
$ws.Range("K105").Value2 = 1241
$ws.Range("M105").Value2 = 506

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 71468230
$ws.Range("J4").Value2 = 10589.8
$ws.Range("L4").Value2 = 31769.4
$ws.Range("N4").Value2 = -31993.4
$ws.Range("H12").Value2 = 258.68182
$ws.Range("J12").Value2 = 277.26666
$ws.Range("L12").Value2 = 831.79998
$ws.Range("N12").Value2 = -1177.79998
$ws.Range("H46").Value2 = 101999.9
$ws.Range("I46").Value2 = 1624.875
$ws.Range("K46").Value2 = 4874.625
$ws.Range("M46").Value2 = -4783.625
$ws.Range("H68").Value2 = 518
$ws.Range("I68").Value2 = 425.16666
$ws.Range("K68").Value2 = 1275.49998
$ws.Range("M68").Value2 = -464.4999800000001
$ws.Range("H71").Value2 = 518
$ws.Range("I71").Value2 = 425.16666
$ws.Range("K71").Value2 = 3826.49994
$ws.Range("M71").Value2 = 229.5000600000003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value2 = 26446.5
$ws.Range("J101").Value2 = 26446.5
$ws.Range("L101").Value2 = 26446.5
$ws.Range("N101").Value2 = -32936.5
$ws.Range("H103").Value2 = 0
$ws.Range("J103").Value2 = 0
$ws.Range("L103").Value2 = 0
$ws.Range("N103").ClearContents()
$ws.Range("H122").Value2 = 4723.75
$ws.Range("I122").Value2 = 3712.3572
$ws.Range("J122").Value2 = 7083.6665
$ws.Range("K122").Value2 = 11137.0716
$ws.Range("L122").Value2 = 21250.9995
$ws.Range("M122").Value2 = -8687.071599999999
$ws.Range("N122").Value2 = -26150.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 1921.6666
$ws.Range("I93").Value2 = 1286.5
$ws.Range("K93").Value2 = 1286.5
$ws.Range("M93").Value2 = -38.5
$ws.Range("H138").Value2 = 106820.75
$ws.Range("J138").Value2 = 106820.75
$ws.Range("L138").Value2 = 106820.75
$ws.Range("N138").Value2 = -117100.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value2 = 11668
$ws.Range("I14").Value2 = 11502
$ws.Range("J14").Value2 = 12000
$ws.Range("K14").Value2 = 11502
$ws.Range("L14").Value2 = 12000
$ws.Range("M14").Value2 = -11334
$ws.Range("N14").Value2 = -12336
$ws.Range("H45").Value2 = 28628.727
$ws.Range("J45").Value2 = 18474.625
$ws.Range("L45").Value2 = 18474.625
$ws.Range("N45").Value2 = -19456.625
$ws.Range("H54").Value2 = 0
$ws.Range("I54").Value2 = 0
$ws.Range("K54").Value2 = 0
$ws.Range("M54").ClearContents()
$ws.Range("H62").Value2 = 6151.143
$ws.Range("I62").Value2 = 6151.143
$ws.Range("K62").Value2 = 6151.143
$ws.Range("M62").Value2 = -5527.143
$ws.Range("H65").Value2 = 6151.143
$ws.Range("I65").Value2 = 6151.143
$ws.Range("K65").Value2 = 30755.715
$ws.Range("M65").Value2 = -27635.715
$ws.Range("H97").Value2 = 46397.5
$ws.Range("J97").Value2 = 46397.5
$ws.Range("L97").Value2 = 46397.5
$ws.Range("N97").Value2 = -48379.5
$ws.Range("H105").Value2 = 36749.5
$ws.Range("J105").Value2 = 36749.5
$ws.Range("L105").Value2 = 36749.5
$ws.Range("N105").Value2 = -43737.5
$ws.Range("H114").Value2 = 125000
$ws.Range("J114").Value2 = 125000
$ws.Range("L114").Value2 = 125000
$ws.Range("N114").Value2 = -133678
$ws.Range("H119").Value2 = 1542424.8
$ws.Range("J119").Value2 = 56566.332
$ws.Range("L119").Value2 = 56566.332
$ws.Range("N119").Value2 = -66242.33199999999
$ws.Range("H122").Value2 = 6518.0586
$ws.Range("I122").Value2 = 4619.1816
$ws.Range("K122").Value2 = 13857.5448
$ws.Range("M122").Value2 = -11407.5448
